$p = $ppt.ActivePresentation

$oldDate = "05-Jan-17"
$newDate = "17-Aug-17"

function Update-DatePlaceholder($shapes) {
    $n = $shapes.Placeholders.Count
    for ($i = 1; $i -le $n; $i++) {
        $ph = $shapes.Placeholders.Item($i)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            $tr = $ph.TextFrame.TextRange
            $found = $tr.Find($oldDate)
            if ($found -ne $null) {
                $found.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# Slide 4: merge the "be " / "established " runs into a single run.
$slide4 = $p.Slides.Item(4)
for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
    $shp = $slide4.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*be established*") {
            $found = $tr.Find("be established ")
            if ($found -ne $null) {
                $found.Text = "be established "
            }
        }
    }
}
